$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Rename the column headers used for the battery alarm/standby detail columns
$ws.Range("L8").Value = "Alarm Current(A)"
$ws.Range("M8").Value = "Standby Current(A)"
$ws.Range("L9").Value = "Alarm Current(A)"
$ws.Range("M9").Value = "Standby Current(A)"
$ws.Range("L10").Value = "Alarm Current(A)"
$ws.Range("M10").Value = "Standby Current(A)"
$ws.Range("L11").Value = "Alarm Current(A)"
$ws.Range("M11").Value = "Standby Current(A)"
